# Update "想去人数" (number of people interested) counts that changed
# between data refreshes for two events, on both the "展览" sheet and
# the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5297
    $ws.Range("F4").Value = 919
}
